$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the source data which stores these as text strings with
# thousand-separator dots), otherwise Excel will auto-convert them to numbers.
$textForceCells = @("D5", "D6", "D11", "D14", "D19", "D20", "D21", "D23", "D28", "D32", "D34", "D35", "D38", "D39", "D42", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.968.58"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "3.420.40"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "577.95"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "152.63"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").Value = "  +4.51%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").Value = "0.418"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "4.005.71"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "28.77"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "3.424.28"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "62.064.67"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "14.52"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "8.96"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("D21").Value = "382.50"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "75.27"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "3.562.94"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "7.72"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "23.22"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "5.54"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").Value = "31.16"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").Value = "168.33"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "3.453.74"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").Value = "42.77"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "2.550.66"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "22.58"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("E51").Value = "  +0.14%  "

# Restore default (General) style on the text-forced cells so no stray
# number-format styling is left behind on cells that should look unchanged.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
